# Generate Report for Handback
# Row 6 on the zh-cn / de-de sheets corresponds to file
# 5fbb0e10-e704-4a68-9b7e-54905c2ec4e2.md. The handback for that file has
# now been processed: we record the (stale) handback xliff as the latest
# target/handback file, stamp the handback datetime, flag that the
# handback commit isn't the latest one, and widen the "Error Detail"
# column so the message is readable.

$wb = $excel.ActiveWorkbook

$fileBase = "5fbb0e10-e704-4a68-9b7e-54905c2ec4e2"
$mdName = "$fileBase.md"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2d80097211e6dd89f47ecfb356c18ae6dffd641/e2e/$mdName, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1655d7a11d286fd1bcfe4f9766e9c3c6d986fa44/e2e/$mdName."

$sheetsInfo = @(
    @{ Name = "zh-cn"; Lang = "zh-cn"; Commit = "ebbc60441908b57e33e051f80ec33092d52c2260"; Repo = "ol-test0-zhcn"; HandbackTime = "2016-10-25 02:37:48" },
    @{ Name = "de-de"; Lang = "de-de"; Commit = "52932682d320f0b06ffac91c8db236117edaea0b"; Repo = "ol-test0-dede"; HandbackTime = "2016-10-25 02:38:05" }
)

foreach ($info in $sheetsInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Widen the "Error Detail" column (16th / P) so the long message fits.
    $ws.Columns.Item(16).ColumnWidth = 39.17

    $xlfName = "$fileBase.40fb1cb65c63877083b8cfe63e3aebf2cdfce4e6.$($info.Lang).xlf"

    # I6: Latest Target File -> hyperlink to the source .md file.
    $ws.Range("I6").Value = $mdName
    $targetUrl = "https://github.com/OpenLocalizationTestOrg/$($info.Repo)/blob/$($info.Commit)/e2e/$mdName"
    $ws.Hyperlinks.Add($ws.Range("I6"), $targetUrl, [System.Type]::Missing, [System.Type]::Missing, $mdName)

    # J6: Latest Handback File -> the xliff that was handed back.
    $ws.Range("J6").Value = $xlfName

    # K6: Latest Handback DateTime.
    $ws.Range("K6").Value = $info.HandbackTime

    # P6: Error Detail -> version mismatch warning.
    $ws.Range("P6").Value = $errorDetail
}
